$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, "D").Value = '''34.702.33'
$ws.Cells.Item(2, "E").Value = '  -0.66%  '

$ws.Cells.Item(3, "D").Value = '''1.820.08'
$ws.Cells.Item(3, "E").Value = '  +0.48%  '

$ws.Cells.Item(4, "E").Value = '  +0.41%  '

$ws.Cells.Item(5, "D").Value = '''230.16'
$ws.Cells.Item(5, "E").Value = '  -1.03%  '

$ws.Cells.Item(6, "E").Value = '  +0.91%  '

$ws.Cells.Item(7, "E").Value = '  +0.39%  '

$ws.Cells.Item(8, "D").Value = '''39.27'
$ws.Cells.Item(8, "E").Value = '  -2.14%  '

$ws.Cells.Item(9, "D").Value = '''0.322'
$ws.Cells.Item(9, "E").Value = '  +1.78%  '

$ws.Cells.Item(10, "E").Value = '  -0.55%  '

$ws.Cells.Item(11, "E").Value = '  -1.00%  '

$ws.Cells.Item(12, "D").Value = '''2.084.68'
$ws.Cells.Item(12, "E").Value = '  +0.54%  '

$ws.Cells.Item(13, "B").Value = 'Chainlink'
$ws.Cells.Item(13, "C").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(13, "D").Value = '''11.23'
$ws.Cells.Item(13, "E").Value = '  +1.26%  '

$ws.Cells.Item(14, "B").Value = 'WrappedEther'
$ws.Cells.Item(14, "C").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, "D").Value = '''1.805.97'
$ws.Cells.Item(14, "E").Value = '  -0.10%  '

$ws.Cells.Item(15, "E").Value = '  +0.91%  '

$ws.Cells.Item(16, "E").Value = '  -1.03%  '

$ws.Cells.Item(17, "D").Value = '''34.602.85'
$ws.Cells.Item(17, "E").Value = '  -0.82%  '

$ws.Cells.Item(18, "D").Value = '''69.29'
$ws.Cells.Item(18, "E").Value = '  +0.33%  '

$ws.Cells.Item(19, "E").Value = '  -0.62%  '

$ws.Cells.Item(20, "D").Value = '''238.42'
$ws.Cells.Item(20, "E").Value = '  +0.65%  '

$ws.Cells.Item(21, "D").Value = '''12.05'
$ws.Cells.Item(21, "E").Value = '  +1.84%  '

$ws.Cells.Item(22, "D").Value = '''4.63'

$ws.Cells.Item(23, "E").Value = '  +0.41%  '

$ws.Cells.Item(24, "D").Value = '''2.25'
$ws.Cells.Item(24, "E").Value = '  -1.34%  '

$ws.Cells.Item(25, "D").Value = '''172.63'
$ws.Cells.Item(25, "E").Value = '  +0.01%  '

$ws.Cells.Item(26, "E").Value = '  -1.40%  '

$ws.Cells.Item(27, "E").Value = '  +3.26%  '

$ws.Cells.Item(28, "D").Value = '''17.26'
$ws.Cells.Item(28, "E").Value = '  -0.65%  '

$ws.Cells.Item(29, "D").Value = '''1.50'
$ws.Cells.Item(29, "E").Value = '  -8.43%  '

$ws.Cells.Item(30, "E").Value = '  +0.33%  '

$ws.Cells.Item(31, "D").Value = '''0.0547'
$ws.Cells.Item(31, "E").Value = '  -1.01%  '

$ws.Cells.Item(32, "D").Value = '''3.88'
$ws.Cells.Item(32, "E").Value = '  +0.12%  '

$ws.Cells.Item(33, "E").Value = '  -1.82%  '

$ws.Cells.Item(34, "E").Value = '  +8.22%  '

$ws.Cells.Item(35, "E").Value = '  +2.41%  '

$ws.Cells.Item(36, "E").Value = '  +11.35%  '

$ws.Cells.Item(37, "E").Value = '  +2.59%  '

$ws.Cells.Item(38, "D").Value = '''91.11'
$ws.Cells.Item(38, "E").Value = '  -1.92%  '

$ws.Cells.Item(39, "D").Value = '''1.337.23'
$ws.Cells.Item(39, "E").Value = '  +2.62%  '

$ws.Cells.Item(40, "E").Value = '  +2.64%  '

$ws.Cells.Item(41, "E").Value = '  +0.00%  '

$ws.Cells.Item(42, "D").Value = '''14.39'
$ws.Cells.Item(42, "E").Value = '  -2.21%  '

$ws.Cells.Item(43, "D").Value = '''2.42'
$ws.Cells.Item(43, "E").Value = '  -0.73%  '

$ws.Cells.Item(44, "E").Value = '  -4.17%  '

$ws.Cells.Item(45, "E").Value = '  -0.75%  '

$ws.Cells.Item(46, "B").Value = 'FraxShare'
$ws.Cells.Item(46, "C").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(46, "D").Value = '''6.26'
$ws.Cells.Item(46, "E").Value = '  +0.43%  '

$ws.Cells.Item(47, "B").Value = 'Kaspa'
$ws.Cells.Item(47, "C").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(47, "D").Value = '''0.0521'
$ws.Cells.Item(47, "E").Value = '  +1.67%  '

$ws.Cells.Item(48, "D").Value = '''1.999.42'
$ws.Cells.Item(48, "E").Value = '  +0.50%  '

$ws.Cells.Item(49, "E").Value = '  +0.37%  '

$ws.Cells.Item(50, "D").Value = '''0.0668'
$ws.Cells.Item(50, "E").Value = '  +3.70%  '

$ws.Cells.Item(51, "B").Value = 'THORChain'
$ws.Cells.Item(51, "C").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(51, "D").Value = '''3.19'
$ws.Cells.Item(51, "E").Value = '  +13.08%  '

